# "Cambios en la distribucion de clases."
#
# Two functional edits to the worksheet data/UI state:
#   1. Fix the decimal separator typo in Ana Torres Pardo's "Localizacion"
#      value: "21.26; 50,26" -> "21.26; 50.26" (cell B4).
#   2. Move the active cell/selection from D4 to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the coordinate string stored in B4 (comma -> period).
$ws.Range("B4").Value = "21.26; 50.26"

# 2) Update the sheet's active selection to B5.
$ws.Range("B5").Select()
